$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting the existing rows 7-9 down to 8-10.
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new result entry.
$ws.Range("A7").Value = "Torneo FEG"
$ws.Range("B7").Value = "Principiantes"
$ws.Range("C7").Value = "general"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "Vera, Bautista"
$ws.Range("F7").Value = 34
